$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on cells whose new value would otherwise be
# auto-converted to a numeric type by Excel (values that look like plain numbers).
foreach ($addr in @("D5","D6","D8","D9","D11","D14","D15","D17","D19","D20","D23","D25","D26","D27","D31","D32","D35","D36","D37","D39","D40","D41","D44","D45","D48","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values.
$ws.Range('D2').Value = '26.774.79'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.538.96'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '205.80'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = '0.481'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '0.245'
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').Value = '21.28'
$ws.Range('E9').Value = '  -3.13%  '
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').Value = '0.0853'
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').Value = '1.759.44'
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('D13').Value = '1.537.94'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').Value = '3.67'
$ws.Range('E14').Value = '  -1.98%  '
$ws.Range('D15').Value = '0.506'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').Value = '26.762.84'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '61.11'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '0.0₃0688'
$ws.Range('E18').Value = '  +1.40%  '
$ws.Range('D19').Value = '211.92'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('D20').Value = '7.23'
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('E22').Value = '  -1.37%  '
$ws.Range('D23').Value = '9.04'
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('E24').Value = '  -0.67%  '
$ws.Range('D25').Value = '152.90'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').Value = '6.47'
$ws.Range('E26').Value = '  -3.67%  '
$ws.Range('D27').Value = '14.83'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('E30').Value = '  -0.78%  '
$ws.Range('D31').Value = '0.0455'
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('D32').Value = '3.22'
$ws.Range('E32').Value = '  +2.05%  '
$ws.Range('D33').Value = '1.360.86'
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').Value = '1.51'
$ws.Range('E35').Value = '  -3.14%  '
$ws.Range('D36').Value = '2.27'
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('D37').Value = '0.929'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('D39').Value = '0.522'
$ws.Range('E39').Value = '  +1.89%  '
$ws.Range('D40').Value = '5.71'
$ws.Range('E40').Value = '  +5.07%  '
$ws.Range('D41').Value = '0.798'
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('D44').Value = '1.74'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('D45').Value = '62.50'
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('D46').Value = '1.673.15'
$ws.Range('E46').Value = '  -1.53%  '
$ws.Range('E47').Value = '  -4.12%  '
$ws.Range('D48').Value = '85.75'
$ws.Range('E49').Value = '  +3.33%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0974'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0949'
$ws.Range('E51').Value = '  +0.43%  '
